$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value2 = 1304.5555
$ws.Range("I41").Value2 = 1304.5555
$ws.Range("K41").Value2 = 1304.5555
$ws.Range("M41").Value2 = -864.5554999999999

$ws.Range("H113").Value2 = 7599.8184
$ws.Range("I113").Value2 = 8128.4287
$ws.Range("K113").Value2 = 8128.4287
$ws.Range("M113").Value2 = -4874.4287

$ws.Range("H137").Value2 = 37038860
$ws.Range("I137").Value2 = 41667904
$ws.Range("J137").Value2 = 6500
$ws.Range("K137").Value2 = 125003712
$ws.Range("L137").Value2 = 19500
$ws.Range("M137").Value2 = -125001162
$ws.Range("N137").Value2 = -24600

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value2 = 1578.25
$ws.Range("I61").Value2 = 1578.25
$ws.Range("K61").Value2 = 1578.25
$ws.Range("M61").Value2 = -1366.25

$ws.Range("H132").Value2 = 9320.305
$ws.Range("I132").Value2 = 9415.637000000001
$ws.Range("J132").Value2 = 9232.916999999999
$ws.Range("K132").Value2 = 28246.911
$ws.Range("L132").Value2 = 27698.751
$ws.Range("M132").Value2 = -25716.911
$ws.Range("N132").Value2 = -32758.751

$ws.Range("H136").Value2 = 1578.25
$ws.Range("I136").Value2 = 1578.25
$ws.Range("K136").Value2 = 4734.75
$ws.Range("M136").Value2 = -2184.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value2 = 7847.5
$ws.Range("I97").Value2 = 7847.5
$ws.Range("K97").Value2 = 7847.5
$ws.Range("M97").Value2 = -6856.5

$ws.Range("H134").Value2 = 1875.125
$ws.Range("I134").Value2 = 1875.125
$ws.Range("K134").Value2 = 5625.375
$ws.Range("M134").Value2 = -3090.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 7110.8125
$ws.Range("I31").Value2 = 9199.111000000001
$ws.Range("J31").Value2 = 4425.857
$ws.Range("K31").Value2 = 9199.111000000001
$ws.Range("L31").Value2 = 4425.857
$ws.Range("M31").Value2 = -8904.111000000001
$ws.Range("N31").Value2 = -5015.857

$ws.Range("H34").Value2 = 7110.8125
$ws.Range("I34").Value2 = 9199.111000000001
$ws.Range("J34").Value2 = 4425.857
$ws.Range("K34").Value2 = 9199.111000000001
$ws.Range("L34").Value2 = 4425.857
$ws.Range("M34").Value2 = -8997.111000000001
$ws.Range("N34").Value2 = -4829.857

$ws.Range("H132").Value2 = 15388580
$ws.Range("I132").Value2 = 22225308
$ws.Range("J132").Value2 = 5941.25
$ws.Range("K132").Value2 = 66675924
$ws.Range("L132").Value2 = 17823.75
$ws.Range("M132").Value2 = -66673394
$ws.Range("N132").Value2 = -22883.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value2 = 1788.9375
$ws.Range("J131").Value2 = 1741.8937
$ws.Range("L131").Value2 = 5225.6811
$ws.Range("N131").Value2 = -15305.6811

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value2 = 5407.6665
$ws.Range("I70").Value2 = 5519.4
$ws.Range("J70").Value2 = 4849
$ws.Range("K70").Value2 = 5519.4
$ws.Range("L70").Value2 = 4849
$ws.Range("M70").Value2 = -5249.4
$ws.Range("N70").Value2 = -5389

$ws.Range("H73").Value2 = 5407.6665
$ws.Range("I73").Value2 = 5519.4
$ws.Range("J73").Value2 = 4849
$ws.Range("K73").Value2 = 5519.4
$ws.Range("L73").Value2 = 4849
$ws.Range("M73").Value2 = -4583.4
$ws.Range("N73").Value2 = -6721

$ws.Range("H80").Value2 = 18531.637
$ws.Range("J80").Value2 = 18399.8
$ws.Range("L80").Value2 = 18399.8
$ws.Range("N80").Value2 = -20395.8

$ws.Range("H83").Value2 = 18531.637
$ws.Range("J83").Value2 = 18399.8
$ws.Range("L83").Value2 = 91999
$ws.Range("N83").Value2 = -101983

$ws.Range("H107").Value2 = 2646
$ws.Range("I107").Value2 = 748.6667
$ws.Range("J107").Value2 = 5492
$ws.Range("K107").Value2 = 748.6667
$ws.Range("L107").Value2 = 5492
$ws.Range("M107").Value2 = 1171.3333
$ws.Range("N107").Value2 = -9332

$ws.Range("H132").Value2 = 41670668
$ws.Range("I132").Value2 = 3497
$ws.Range("J132").Value2 = 55559724
$ws.Range("K132").Value2 = 10491
$ws.Range("L132").Value2 = 166679172
$ws.Range("M132").Value2 = -7961
$ws.Range("N132").Value2 = -166684232

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 3809.4285
$ws.Range("I7").Value2 = 3815.6
$ws.Range("K7").Value2 = 3815.6
$ws.Range("M7").Value2 = -3703.6

$ws.Range("H22").Value2 = 55556760
$ws.Range("I22").Value2 = 975
$ws.Range("J22").Value2 = 71429840
$ws.Range("K22").Value2 = 975
$ws.Range("L22").Value2 = 71429840
$ws.Range("M22").Value2 = -680
$ws.Range("N22").Value2 = -71430430

$ws.Range("H27").Value2 = 55556760
$ws.Range("I27").Value2 = 975
$ws.Range("J27").Value2 = 71429840
$ws.Range("K27").Value2 = 975
$ws.Range("L27").Value2 = 71429840
$ws.Range("M27").Value2 = -868
$ws.Range("N27").Value2 = -71430054

$ws.Range("H40").Value2 = 1000000000
$ws.Range("I40").Value2 = 1000000000
$ws.Range("J40").Value2 = 0
$ws.Range("K40").Value2 = 1000000000
$ws.Range("L40").Value2 = 0
$ws.Range("M40").Value2 = -999999864
$ws.Range("N40").ClearContents()

$ws.Range("H55").Value2 = 875.43475
$ws.Range("I55").Value2 = 422.15384
$ws.Range("K55").Value2 = 422.15384
$ws.Range("M55").Value2 = -249.15384

$ws.Range("H82").Value2 = 1501.6666
$ws.Range("I82").Value2 = 0
$ws.Range("J82").Value2 = 1501.6666
$ws.Range("K82").Value2 = 0
$ws.Range("L82").Value2 = 1501.6666
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value2 = -2223.6666

$ws.Range("H85").Value2 = 1501.6666
$ws.Range("I85").Value2 = 0
$ws.Range("J85").Value2 = 1501.6666
$ws.Range("K85").Value2 = 0
$ws.Range("L85").Value2 = 1501.6666
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value2 = -3997.6666

$ws.Range("H101").Value2 = 6025.7
$ws.Range("J101").Value2 = 6025.7
$ws.Range("L101").Value2 = 6025.7
$ws.Range("N101").Value2 = -12515.7

$ws.Range("H126").Value2 = 3809.4285
$ws.Range("I126").Value2 = 3815.6
$ws.Range("K126").Value2 = 11446.8
$ws.Range("M126").Value2 = -8976.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value2 = 35000
$ws.Range("J51").Value2 = 40000
$ws.Range("L51").Value2 = 40000
$ws.Range("N51").Value2 = -41020

$ws.Range("H52").Value2 = 15000
$ws.Range("I52").Value2 = 15000
$ws.Range("K52").Value2 = 15000
$ws.Range("M52").Value2 = -14774

$ws.Range("H96").Value2 = 2874.5
$ws.Range("I96").Value2 = 4599
$ws.Range("J96").Value2 = 1150
$ws.Range("K96").Value2 = 4599
$ws.Range("L96").Value2 = 1150
$ws.Range("M96").Value2 = -3226
$ws.Range("N96").Value2 = -3896

$ws.Range("H132").Value2 = 250019000
$ws.Range("I132").Value2 = 35004
$ws.Range("J132").Value2 = 500003000
$ws.Range("K132").Value2 = 105012
$ws.Range("L132").Value2 = 1500009000
$ws.Range("M132").Value2 = -102482
$ws.Range("N132").Value2 = -1500014060

$ws.Range("H136").Value2 = 55888.11
$ws.Range("I136").Value2 = 56124.125
$ws.Range("K136").Value2 = 168372.375
$ws.Range("M136").Value2 = -165822.375
